$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.815.48"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "1.643.80"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'216.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").Value = "'0.503"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "'1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").Value = "'0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("E9").Value = "  -0.49%  "
$ws.Range("D10").Value = "'19.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.21%  "
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.873.67"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").Value = "1.645.81"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").Value = "'4.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.30%  "
$ws.Range("D15").Value = "'0.526"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("D16").Value = "'66.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.80%  "
$ws.Range("D17").Value = "26.844.67"
$ws.Range("E17").Value = "  +0.56%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'218.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.19%  "
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  +7.78%  "
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "'2.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.04%  "
$ws.Range("D24").Value = "'9.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").Value = "'145.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.49%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("E27").Value = "  +3.58%  "
$ws.Range("E28").Value = "  +0.32%  "
$ws.Range("D29").Value = "'15.83"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.75%  "
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'3.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("E34").Value = "  +1.33%  "
$ws.Range("D35").Value = "'2.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Value = "1.244.84"
$ws.Range("E36").Value = "  -2.17%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "'0.532"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.88%  "
$ws.Range("E39").Value = "  +2.83%  "
$ws.Range("D40").Value = "'1.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "'0.805"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("D42").Value = "'5.35"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("D43").Value = "1.787.06"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "'2.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.86%  "
$ws.Range("D45").Value = "'60.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'91.75"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.33%  "
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("D50").Value = "'0.0970"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("E51").Value = "  -0.09%  "
